$d = $word.ActiveDocument

# Remove the "_GoBack" bookmark (bookmarkStart/bookmarkEnd wrapping the title run)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# The document currently ends with an empty paragraph (3rd paragraph).
# Replace its (empty) text with "change from rhythm", then append a new
# paragraph "30th of July" where "th" is a superscript.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Text = "change from rhythm"

# Move to the end of the document and insert a new paragraph with the date.
$endRange = $d.Content
$endRange.Collapse(0)  # wdCollapseEnd
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertAfter("30")

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertAfter("th")
$thRange = $d.Content
$thRange.Collapse(0)
$thStart = $thRange.Start - 2
$thEnd = $thRange.Start
$thRangeObj = $d.Range($thStart, $thEnd)
$thRangeObj.Font.Superscript = $true

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertAfter(" of July")
